$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of Oman M2 data to append (rows 291-293)
$newRows = @(
    @{ Row = 291; Date = 44958.45833333334; Symbol = "ECONOMICS:OMM2"; Open = 20665600000; High = 20665600000; Low = 20665600000; Close = 20665600000; Volume = 0 },
    @{ Row = 292; Date = 44986.45833333334; Symbol = "ECONOMICS:OMM2"; Open = 21052600000; High = 21052600000; Low = 21052600000; Close = 21052600000; Volume = 0 },
    @{ Row = 293; Date = 45017.45833333334; Symbol = "ECONOMICS:OMM2"; Open = 21069600000; High = 21069600000; Low = 21069600000; Close = 21069600000; Volume = 0 }
)

$lastExistingRow = 290

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Copy formatting from the last existing data row so the new rows match
    $ws.Range("A" + $lastExistingRow + ":G" + $lastExistingRow).Copy()
    $ws.Range("A" + $rowIndex + ":G" + $rowIndex).PasteSpecial(-4122)

    # Column A: datetime value styled like the preceding rows
    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.Value = $r.Date

    $ws.Cells.Item($rowIndex, 2).Value = $r.Symbol
    $ws.Cells.Item($rowIndex, 3).Value = $r.Open
    $ws.Cells.Item($rowIndex, 4).Value = $r.High
    $ws.Cells.Item($rowIndex, 5).Value = $r.Low
    $ws.Cells.Item($rowIndex, 6).Value = $r.Close
    $ws.Cells.Item($rowIndex, 7).Value = $r.Volume
}

$excel.CutCopyMode = 0

$wb.Save()
